# ----------------------------------------------------------------------
# Rebuild the "module breakdown" sheet:
#  - insert a new column before column B (old B..J shift to C..K)
#  - add a sub-feature/feature breakdown table in rows 3-12
#  - give every cell in the table a thin box border
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column B -> old columns B..J become C..K
$ws.Columns("B:B").Insert()

# 2) Values, entered in the same order as the original authoring session
#    (keeps shared-string table ordering identical to the source file).
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

$ws.Range("C3").Value = "Dashboard"
$ws.Range("C4").Value = "Letter Sent"
$ws.Range("C5").Value = "Exclusion"
$ws.Range("C6").Value = "Exception"
$ws.Range("C7").Value = "Special Request"
$ws.Range("B4").Value = "Letter Search"
$ws.Range("B8").Value = "Files"
$ws.Range("C8").Value = "Upload"
$ws.Range("C9").Value = "Download"
$ws.Range("B10").Value = "Maintenance"
$ws.Range("C10").Value = "Delete Letter"

# 3) Formatting pass - every individual cell in the table gets a thin box
#    border (done before merging so the merge-group cells start from the
#    same per-cell box look).
$ws.Range("A1:K12").Borders.LineStyle = 1

# Header rows 1-2 (B:K) are bold; A1/A2 drop back to regular weight.
$ws.Range("B1:K2").Font.Bold = $true
$ws.Range("A1:A2").Font.Bold = $false

# Wrap the long sub-header label.
$ws.Range("D2").WrapText = $true

# 4) Merges.
# Row-1 category headers: already bold/boxed/centred per-cell above, so a
# plain Merge() is enough - every constituent cell keeps its full box
# border (that's how the source file looks: D1/E1/F1 all carry the same
# bold+border+center style).
$ws.Range("D1:F1").Merge()
$ws.Range("G1:I1").Merge()
$ws.Range("J1:K1").Merge()
$ws.Range("D1:F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("G1:I1").HorizontalAlignment = -4108
$ws.Range("J1:K1").HorizontalAlignment = -4108

# Grouping cells in column B: re-applying the 4 outer edges after merging
# turns the per-cell boxes of the now-merged range into a single outline
# (top edge only on the first row, bottom edge only on the last row) -
# exactly how Excel redraws a bordered block once it becomes one cell.
$ws.Range("B4:B7").Merge()
$ws.Range("B4:B7").HorizontalAlignment = -4108
$ws.Range("B4:B7").Borders.Item(7).LineStyle = 1
$ws.Range("B4:B7").Borders.Item(10).LineStyle = 1
$ws.Range("B4:B7").Borders.Item(8).LineStyle = 1
$ws.Range("B4:B7").Borders.Item(9).LineStyle = 1

$ws.Range("B8:B9").Merge()
$ws.Range("B8:B9").HorizontalAlignment = -4108
$ws.Range("B8:B9").Borders.Item(7).LineStyle = 1
$ws.Range("B8:B9").Borders.Item(10).LineStyle = 1
$ws.Range("B8:B9").Borders.Item(8).LineStyle = 1
$ws.Range("B8:B9").Borders.Item(9).LineStyle = 1

# 5) Selection parity with the saved workbook
$ws.Range("D6").Select()
